# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# values, as produced by the scheduled "Updated cryptos list" GitHub Action.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to keep its text representation (e.g. trailing
# zeros such as "69.20") instead of Excel auto-coercing numeric-looking
# strings into Number cells.
$priceRows = 2,3,5,6,8,12,13,14,15,16,17,18,19,20,21,22,24,25,26,28,30,32,33,34,38,40,41,42,43,46,48,49,51
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "37.318.12"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "2.063.50"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "233.44"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "56.55"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Value = "2.363.40"
$ws.Range("E12").Value = "  -0.82%  "

$ws.Range("D13").Value = "14.63"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").Value = "20.61"
$ws.Range("E14").Value = "  -2.87%  "

$ws.Range("D15").Value = "0.776"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").Value = "5.12"
$ws.Range("E16").Value = "  -2.53%  "

$ws.Range("D17").Value = "2.058.19"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").Value = "37.248.35"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").Value = "6.32"
$ws.Range("E19").Value = "  +6.18%  "

$ws.Range("D20").Value = "69.20"
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("D21").Value = "0.0₃0808"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").Value = "225.93"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  -2.46%  "

$ws.Range("D26").Value = "166.55"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("E27").Value = "  +4.98%  "

$ws.Range("D28").Value = "8.75"
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("E29").Value = "  -3.42%  "

$ws.Range("D30").Value = "19.02"
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("D32").Value = "4.46"
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").Value = "0.0615"
$ws.Range("E33").Value = "  -1.00%  "

$ws.Range("D34").Value = "4.55"
$ws.Range("E34").Value = "  +3.85%  "

$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -1.44%  "

$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  -3.40%  "

$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").Value = "1.465.19"
$ws.Range("E41").Value = "  -0.49%  "

$ws.Range("D42").Value = "96.15"
$ws.Range("E42").Value = "  +1.24%  "

$ws.Range("D43").Value = "0.0932"
$ws.Range("E43").Value = "  -2.58%  "

$ws.Range("E44").Value = "  +1.56%  "

$ws.Range("E45").Value = "  +2.57%  "

$ws.Range("D46").Value = "4.23"
$ws.Range("E46").Value = "  -5.87%  "

$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").Value = "15.07"
$ws.Range("E48").Value = "  -6.75%  "

$ws.Range("D49").Value = "7.15"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("D51").Value = "2.249.99"

